# fix: corrected wrong dates in journal and summary tables
#
# Row 4 (date 2024-01-25) note text: "kvallen 23-01-23" -> "kvallen 23-01-24"
# Row 6 (date 2024-01-27) note text: "utfordes 24-01-24" -> "utfordes 24-01-25"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$e4 = $ws.Range("E4").Value2
$e4Fixed = $e4.Replace("23-01-23", "23-01-24")
$ws.Range("E4").Value = $e4Fixed
$ws.Rows(4).RowHeight = 14.25

$e6 = $ws.Range("E6").Value2
$e6Fixed = $e6.Replace("24-01-24", "24-01-25")
$ws.Range("E6").Value = $e6Fixed
$ws.Rows(6).RowHeight = 14.25

$ws.Range("E9").Select() | Out-Null
